# Generate Report for Handback
# Update the "Correspond Handoff Datetime" (E3) and "Correspond Handback
# DateTime" (H3) values on the zh-cn and de-de sheets to reflect the new
# report generation timestamps.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E3").Value = "2016-03-25 02:55:15"
$wsZhCn.Range("H3").Value = "2016-03-25 02:55:42"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E3").Value = "2016-03-25 02:55:19"
$wsDeDe.Range("H3").Value = "2016-03-25 02:55:50"
